# Fruta / hortaliza, semanal
# New weekly price-report rows for "Pepino ensalada" (Vega Monumental Concepción):
#  - a new latest observation is inserted as row 38, pushing the existing
#    rows 38..63 down to 39..64
#  - a further new observation is appended as the new last row (65)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38 (shifts rows 38:63 down to 39:64)
$ws.Rows.Item(38).Insert()

$ws.Cells.Item(38, 1).Value2 = 11
$ws.Cells.Item(38, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(38, 3).Value2 = "Bíobío"
$ws.Cells.Item(38, 4).Value2 = 44435
$ws.Cells.Item(38, 5).Value2 = 8
$ws.Cells.Item(38, 6).Value2 = 100112043
$ws.Cells.Item(38, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(38, 8).Value2 = "Sin especificar"
$ws.Cells.Item(38, 9).Value2 = "Primera"
$ws.Cells.Item(38, 10).Value2 = 100
$ws.Cells.Item(38, 11).Value2 = 14000
$ws.Cells.Item(38, 12).Value2 = 15000
$ws.Cells.Item(38, 13).Value2 = 14500
$ws.Cells.Item(38, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(38, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(38, 16).Value2 = 242
$ws.Cells.Item(38, 17).Value2 = 60
$ws.Cells.Item(38, 18).Value2 = "Hortaliza"
$ws.Cells.Item(38, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append a new row at the end (row 65)
$ws.Cells.Item(65, 1).Value2 = 11
$ws.Cells.Item(65, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(65, 3).Value2 = "Bíobío"
$ws.Cells.Item(65, 4).Value2 = 44432
$ws.Cells.Item(65, 5).Value2 = 8
$ws.Cells.Item(65, 6).Value2 = 100112043
$ws.Cells.Item(65, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(65, 8).Value2 = "Sin especificar"
$ws.Cells.Item(65, 9).Value2 = "Primera"
$ws.Cells.Item(65, 10).Value2 = 100
$ws.Cells.Item(65, 11).Value2 = 14000
$ws.Cells.Item(65, 12).Value2 = 15000
$ws.Cells.Item(65, 13).Value2 = 14500
$ws.Cells.Item(65, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(65, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(65, 16).Value2 = 242
$ws.Cells.Item(65, 17).Value2 = 60
$ws.Cells.Item(65, 18).Value2 = "Hortaliza"
$ws.Cells.Item(65, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
